$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record for Albahaca @ Femacal de La Calera needs to be
# inserted as row 115; every existing record from row 115 down to the last
# row (160) shifts down by one row (to 116:161), and the new row 115 gets
# the new record's data.
$ws.Rows.Item(115).Insert()

# All the columns that stay constant for every record of this sheet
# (mercado/region/categoria/etc.) - set explicitly on the freshly inserted,
# still-empty row 115.
$ws.Range("A115").Value = 3
$ws.Range("B115").Value = "Femacal de La Calera"
$ws.Range("C115").Value = "Coquimbo"
$ws.Range("D115").Value = 44636
$ws.Range("E115").Value = 5
$ws.Range("F115").Value = 100112052
$ws.Range("G115").Value = "Albahaca"
$ws.Range("H115").Value = "Sin especificar"
$ws.Range("I115").Value = "Primera"
$ws.Range("J115").Value = 130
$ws.Range("K115").Value = 4500
$ws.Range("L115").Value = 5000
$ws.Range("M115").Value = 4769
$ws.Range("N115").Value = "$/docena de matas"
$ws.Range("O115").Value = "Provincia de Quillota"
$ws.Range("P115").Value = 795
$ws.Range("Q115").Value = 6
$ws.Range("R115").Value = "Hortaliza"

# Match the date column's number format used by the rest of column D.
$ws.Range("D115").NumberFormat = $ws.Range("D116").NumberFormat

Write-Host "Inserted new row 115; sheet now spans" $ws.UsedRange.Address()
